$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update existing CELIA balance (original row 5) before any row insertions shift things.
$ws.Cells.Item(5,3).Value = 60663.67

# 2) Insert new rows from the BOTTOM of the sheet upward so earlier row indices remain valid.

# Insert RAPHAELA row right after ANGELA (original row 108 / before original row 109, GUSTAVO)
$ws.Rows.Item(109).Insert()
$ws.Cells.Item(109,1).Value = "'005366255"
$ws.Cells.Item(109,2).Value = "RAPHAELA"
$ws.Cells.Item(109,3).Value = 25.96

# Insert PATRICIA row before THOMAS (original row 8)
$ws.Rows.Item(8).Insert()
$ws.Cells.Item(8,1).Value = "'005255637"
$ws.Cells.Item(8,2).Value = "PATRICIA"
$ws.Cells.Item(8,3).Value = 20000

# Insert RACHEL row before THIAGO (original row 7)
$ws.Rows.Item(7).Insert()
$ws.Cells.Item(7,1).Value = "'004588677"
$ws.Cells.Item(7,2).Value = "RACHEL"
$ws.Cells.Item(7,3).Value = 29696.52

# Insert LEVI row before HFR (original row 4)
$ws.Rows.Item(4).Insert()
$ws.Cells.Item(4,1).Value = "'005206566"
$ws.Cells.Item(4,2).Value = "LEVI"
$ws.Cells.Item(4,3).Value = 100000
